$d = $word.ActiveDocument

$replacements = @(
    @("10÷8=", "65÷4="),
    @("63÷6=", "46÷9="),
    @("20÷5=", "33÷5="),
    @("96÷6=", "73÷3="),
    @("35÷4=", "80÷8="),
    @("93÷5=", "79÷8="),
    @("30÷5=", "35÷9="),
    @("80÷2=", "32÷5="),
    @("76÷8=", "31÷5="),
    @("43÷8=", "21÷5="),
    @("89÷9=", "99÷4="),
    @("95÷8=", "73÷3="),
    @("47÷2=", "17÷4="),
    @("38÷9=", "59÷8="),
    @("76÷6=", "53÷9="),
    @("77÷5=", "56÷9="),
    @("26÷4=", "90÷8="),
    @("29÷6=", "73÷3="),
    @("84÷3=", "36÷2="),
    @("97÷3=", "16÷4="),
    @("24÷8=", "30÷8="),
    @("85÷7=", "33÷6="),
    @("85÷3=", "78÷4="),
    @("26÷7=", "69÷2="),
    @("81÷9=", "76÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
